# Update LR-pair (ligand-receptor) TPM-derived expression/specificity metrics
# on the active worksheet to reflect newly recomputed TPM values (per commit:
# "update scripts wuth new tpm"). Columns G:J (ligand), M:P (receptor) and
# Q:T (edge) are refreshed with updated averages/totals/specificities for
# data rows 2-10; columns A:F, K:L are unchanged identifiers/counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.973713333333333
$ws.Cells.Item(2, 8).Value = 5.921139999999999
$ws.Cells.Item(2, 9).Value = 0.03018202516692187
$ws.Cells.Item(2, 10).Value = 0.03018202516692186
$ws.Cells.Item(2, 13).Value = 4.331589999999999
$ws.Cells.Item(2, 14).Value = 12.99477
$ws.Cells.Item(2, 15).Value = 0.1478799966101367
$ws.Cells.Item(2, 16).Value = 0.1478799966101367
$ws.Cells.Item(2, 17).Value = 8.549316937533332
$ws.Cells.Item(2, 18).Value = 76.94385243779999
$ws.Cells.Item(2, 19).Value = 0.004463317779371467
$ws.Cells.Item(2, 20).Value = 0.004463317779371466
$ws.Cells.Item(3, 7).Value = 1.973713333333333
$ws.Cells.Item(3, 8).Value = 5.921139999999999
$ws.Cells.Item(3, 9).Value = 0.03018202516692187
$ws.Cells.Item(3, 10).Value = 0.03018202516692186
$ws.Cells.Item(3, 15).Value = 0.5404313285772905
$ws.Cells.Item(3, 16).Value = 0.5404313285772904
$ws.Cells.Item(3, 17).Value = 31.24370311665778
$ws.Cells.Item(3, 18).Value = 281.19332804992
$ws.Cells.Item(3, 19).Value = 0.0163113119601128
$ws.Cells.Item(3, 20).Value = 0.0163113119601128
$ws.Cells.Item(4, 7).Value = 1.973713333333333
$ws.Cells.Item(4, 8).Value = 5.921139999999999
$ws.Cells.Item(4, 9).Value = 0.03018202516692187
$ws.Cells.Item(4, 10).Value = 0.03018202516692186
$ws.Cells.Item(4, 13).Value = 9.129751000000001
$ws.Cells.Item(4, 15).Value = 0.3116886748125729
$ws.Cells.Item(4, 16).Value = 0.3116886748125729
$ws.Cells.Item(4, 17).Value = 18.01951127871333
$ws.Cells.Item(4, 18).Value = 162.17560150842
$ws.Cells.Item(4, 19).Value = 0.009407395427437602
$ws.Cells.Item(4, 20).Value = 0.0094073954274376
$ws.Cells.Item(5, 7).Value = 43.63696533333334
$ws.Cells.Item(5, 9).Value = 0.667296493191563
$ws.Cells.Item(5, 10).Value = 0.6672964931915629
$ws.Cells.Item(5, 13).Value = 4.331589999999999
$ws.Cells.Item(5, 14).Value = 12.99477
$ws.Cells.Item(5, 15).Value = 0.1478799966101367
$ws.Cells.Item(5, 16).Value = 0.1478799966101367
$ws.Cells.Item(5, 17).Value = 189.0174426682133
$ws.Cells.Item(5, 18).Value = 1701.15698401392
$ws.Cells.Item(5, 19).Value = 0.09867980315112447
$ws.Cells.Item(5, 20).Value = 0.09867980315112446
$ws.Cells.Item(6, 7).Value = 43.63696533333334
$ws.Cells.Item(6, 9).Value = 0.667296493191563
$ws.Cells.Item(6, 10).Value = 0.6672964931915629
$ws.Cells.Item(6, 15).Value = 0.5404313285772905
$ws.Cells.Item(6, 16).Value = 0.5404313285772904
$ws.Cells.Item(6, 17).Value = 690.7692048084765
$ws.Cells.Item(6, 18).Value = 6216.922843276288
$ws.Cells.Item(6, 19).Value = 0.3606279303704833
$ws.Cells.Item(6, 20).Value = 0.3606279303704831
$ws.Cells.Item(7, 7).Value = 43.63696533333334
$ws.Cells.Item(7, 9).Value = 0.667296493191563
$ws.Cells.Item(7, 10).Value = 0.6672964931915629
$ws.Cells.Item(7, 13).Value = 9.129751000000001
$ws.Cells.Item(7, 15).Value = 0.3116886748125729
$ws.Cells.Item(7, 16).Value = 0.3116886748125729
$ws.Cells.Item(7, 17).Value = 398.3946278889654
$ws.Cells.Item(7, 18).Value = 3585.551651000688
$ws.Cells.Item(7, 19).Value = 0.2079887596699554
$ws.Cells.Item(7, 20).Value = 0.2079887596699553
$ws.Cells.Item(8, 7).Value = 19.78298933333334
$ws.Cells.Item(8, 8).Value = 59.34896800000001
$ws.Cells.Item(8, 9).Value = 0.3025214816415151
$ws.Cells.Item(8, 10).Value = 0.3025214816415151
$ws.Cells.Item(8, 13).Value = 4.331589999999999
$ws.Cells.Item(8, 14).Value = 12.99477
$ws.Cells.Item(8, 15).Value = 0.1478799966101367
$ws.Cells.Item(8, 16).Value = 0.1478799966101367
$ws.Cells.Item(8, 17).Value = 85.69179876637334
$ws.Cells.Item(8, 18).Value = 771.22618889736
$ws.Cells.Item(8, 19).Value = 0.0447368756796408
$ws.Cells.Item(8, 20).Value = 0.0447368756796408
$ws.Cells.Item(9, 7).Value = 19.78298933333334
$ws.Cells.Item(9, 8).Value = 59.34896800000001
$ws.Cells.Item(9, 9).Value = 0.3025214816415151
$ws.Cells.Item(9, 10).Value = 0.3025214816415151
$ws.Cells.Item(9, 15).Value = 0.5404313285772905
$ws.Cells.Item(9, 16).Value = 0.5404313285772904
$ws.Cells.Item(9, 17).Value = 313.1629274889672
$ws.Cells.Item(9, 19).Value = 0.1634920862466944
$ws.Cells.Item(9, 20).Value = 0.1634920862466944
$ws.Cells.Item(10, 7).Value = 19.78298933333334
$ws.Cells.Item(10, 8).Value = 59.34896800000001
$ws.Cells.Item(10, 9).Value = 0.3025214816415151
$ws.Cells.Item(10, 10).Value = 0.3025214816415151
$ws.Cells.Item(10, 13).Value = 9.129751000000001
$ws.Cells.Item(10, 15).Value = 0.3116886748125729
$ws.Cells.Item(10, 16).Value = 0.3116886748125729
$ws.Cells.Item(10, 17).Value = 180.6137666489894
$ws.Cells.Item(10, 18).Value = 1625.523899840904
$ws.Cells.Item(10, 19).Value = 0.09429251971517995
$ws.Cells.Item(10, 20).Value = 0.09429251971517995
